$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns (I, J), matching the style of the
# existing header row (row 1).
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill data rows 2..39: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 39; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
